$wb = $excel.ActiveWorkbook

# Sheet 1: display/exhibitions (展览) - update "want to go" counts (column F)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 111
$ws1.Range("F4").Value = 419
$ws1.Range("F6").Value = 133
$ws1.Range("F7").Value = 1160
$ws1.Range("F8").Value = 389
$ws1.Range("F14").Value = 790
$ws1.Range("F15").Value = 180
$ws1.Range("F19").Value = 1013
$ws1.Range("F20").Value = 464
$ws1.Range("F21").Value = 265
$ws1.Range("F22").Value = 85
$ws1.Range("F23").Value = 382
$ws1.Range("F25").Value = 43
$ws1.Range("F26").Value = 473

# Sheet 2: performances (演出) - update counts, fold row 14 into row 13, remove row 14
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F3").Value = 21
$ws2.Range("F4").Value = 367
$ws2.Range("F6").Value = 42
$ws2.Range("B13").Value = "'2024.04.28"
$ws2.Range("C13").Value = " 广州·夏川里美 2024 巡回演唱会 出道 25 周年纪念专场"
$ws2.Range("D13").Value = "中山纪念堂 中山纪念堂"
$ws2.Range("E13").Value = "2024.04.28 19:30-04.28 21:30"
$ws2.Range("F13").Value = 13
$ws2.Range("G13").Value = 280
$ws2.Range("H13").Value = "https://show.bilibili.com/platform/detail.html?id=81068"
$ws2.Range("I13").Value = "//i0.hdslb.com/bfs/openplatform/202401/pXznRv8G1705633441713.jpeg"
$ws2.Rows.Item(14).Delete()

# Sheet 3: local life (本地生活) - unchanged

# Sheet 4: all types (全部类型) - update counts, fold row 40 into row 39, remove row 40
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F4").Value = 111
$ws4.Range("F6").Value = 419
$ws4.Range("F8").Value = 133
$ws4.Range("F9").Value = 1160
$ws4.Range("F10").Value = 389
$ws4.Range("F12").Value = 21
$ws4.Range("F14").Value = 367
$ws4.Range("F18").Value = 42
$ws4.Range("F21").Value = 790
$ws4.Range("F22").Value = 180
$ws4.Range("F26").Value = 1013
$ws4.Range("F27").Value = 464
$ws4.Range("F30").Value = 265
$ws4.Range("F31").Value = 85
$ws4.Range("F32").Value = 382
$ws4.Range("F36").Value = 43
$ws4.Range("F38").Value = 473
$ws4.Range("B39").Value = "'2024.04.28"
$ws4.Range("C39").Value = " 广州·夏川里美 2024 巡回演唱会 出道 25 周年纪念专场"
$ws4.Range("D39").Value = "中山纪念堂 中山纪念堂"
$ws4.Range("E39").Value = "2024.04.28 19:30-04.28 21:30"
$ws4.Range("F39").Value = 13
$ws4.Range("G39").Value = 280
$ws4.Range("H39").Value = "https://show.bilibili.com/platform/detail.html?id=81068"
$ws4.Range("I39").Value = "//i0.hdslb.com/bfs/openplatform/202401/pXznRv8G1705633441713.jpeg"
$ws4.Rows.Item(40).Delete()

